$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Add the two new version rows (1.15 and 1.16) right after the existing last row (14)
$ws.Range("A15").Value = "[1.15]"
$ws.Range("B15").Value = "[Printing]
- print on both sides for all documents by default
- add the margins for envelops only
- use Landscape view for envelops only"

$ws.Range("A16").Value = "[1.16]"
$ws.Range("B16").Value = "[Printing]
- customize margins for envelops
- stop using landscape for envelops"

# Copy styling from row 9, an existing "no date yet" entry whose A/B/C
# cell styles (left/top, wrap, left/top) match what the new rows need
$ws.Range("A9:C9").Copy() | Out-Null
$ws.Range("A15:C15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row heights matching new content (60 for 1.15 with 4 lines, 45 for 1.16 with 3 lines)
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 45

# Clear the C column values (no dates for these entries) but keep style
$ws.Range("C15").ClearContents() | Out-Null
$ws.Range("C16").ClearContents() | Out-Null

# Expand the table (ListObject) range to include the new rows
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:C16"))

# Update the sheet dimension/selection to match the new data
$ws.Range("M15").Select() | Out-Null

$wb.Save()
